$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer footer (A59).
$ws.Range("A59").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-26 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) columns for each holding (rows 2-56)
# with the updated values.
$ws.Range("D2").Value2 = 0.01538701972546092
$ws.Range("E2").Value2 = 0.006829854660692858
$ws.Range("D3").Value2 = 0.05201433810148039
$ws.Range("E3").Value2 = 0.02038983740810796
$ws.Range("D4").Value2 = 0.01519819824953399
$ws.Range("E4").Value2 = 0.01219037472597284
$ws.Range("D5").Value2 = 0.009646585297190336
$ws.Range("E5").Value2 = -0.007488702388637747
$ws.Range("D6").Value2 = 0.0154389580865833
$ws.Range("E6").Value2 = -0.01693753000068565
$ws.Range("D7").Value2 = 0.01960199833155799
$ws.Range("E7").Value2 = -0.02000895923547874
$ws.Range("D8").Value2 = 0.003838612316596943
$ws.Range("E8").Value2 = 0.02055533473397309
$ws.Range("D9").Value2 = 0.006229240419217305
$ws.Range("E9").Value2 = 0.001979485333813313
$ws.Range("D10").Value2 = 0.01414638416059341
$ws.Range("E10").Value2 = 0.0002575328354368356
$ws.Range("D11").Value2 = 0.008227310417119531
$ws.Range("E11").Value2 = -0.01203542502460075
$ws.Range("D12").Value2 = 0.01432601989519474
$ws.Range("E12").Value2 = 0.01147626499739185
$ws.Range("D13").Value2 = 0.002971976545159237
$ws.Range("E13").Value2 = 0.01026769343601042
$ws.Range("D14").Value2 = 0.005922156415986539
$ws.Range("E14").Value2 = 0.00809716599190291
$ws.Range("D15").Value2 = 0.01356223328849757
$ws.Range("E15").Value2 = 0.00246354617484501
$ws.Range("D16").Value2 = 0.01000187108880086
$ws.Range("E16").Value2 = 0.01148781171196411
$ws.Range("D17").Value2 = 0.021762266724596
$ws.Range("E17").Value2 = 0.01918743158843594
$ws.Range("D18").Value2 = 0.008632790835905886
$ws.Range("E18").Value2 = -0.001817906377821532
$ws.Range("D19").Value2 = 0.01649273387621376
$ws.Range("E19").Value2 = -0.008458192363460681
$ws.Range("D20").Value2 = 0.01175276679977896
$ws.Range("E20").Value2 = -0.007948283170835069
$ws.Range("D21").Value2 = 0.007235250353954637
$ws.Range("E21").Value2 = -0.01560509554140121
$ws.Range("D22").Value2 = 0.01326109666354881
$ws.Range("E22").Value2 = -0.003944773175542426
$ws.Range("D23").Value2 = 0.01921233607071538
$ws.Range("E23").Value2 = -0.003588307369158161
$ws.Range("D24").Value2 = 0.00942827603470292
$ws.Range("E24").Value2 = 0.0050761421319796
$ws.Range("D25").Value2 = 0.02062513422469899
$ws.Range("E25").Value2 = 0.005181347150259086
$ws.Range("D26").Value2 = 0.01170219851532887
$ws.Range("E26").Value2 = -0.01760169870362094
$ws.Range("D27").Value2 = 0.02149273527742854
$ws.Range("E27").Value2 = 0.02498844603101236
$ws.Range("D28").Value2 = 0.0585545859253082
$ws.Range("E28").Value2 = 0.002977963073257817
$ws.Range("D29").Value2 = 0.02174607490937799
$ws.Range("E29").Value2 = -0.003891891891891985
$ws.Range("D30").Value2 = 0.03103846424867439
$ws.Range("E30").Value2 = 0.01709470304975924
$ws.Range("D31").Value2 = 0.01560302476809039
$ws.Range("E31").Value2 = 0.02606512589529664
$ws.Range("D32").Value2 = 0.01391026390136619
$ws.Range("E32").Value2 = -0.01395032323919698
$ws.Range("D33").Value2 = 0.01928744118284199
$ws.Range("E33").Value2 = 0.02103917238172737
$ws.Range("D34").Value2 = 0.04296927873963906
$ws.Range("E34").Value2 = 0.004347958416125675
$ws.Range("D35").Value2 = 0.01089584611747205
$ws.Range("E35").Value2 = 0.000685871056241405
$ws.Range("D36").Value2 = 0.009710854348363353
$ws.Range("E36").Value2 = 0.0000897827258035111
$ws.Range("D37").Value2 = 0.01069593947574203
$ws.Range("E37").Value2 = 0.0176855895196506
$ws.Range("D38").Value2 = 0.007351706871099543
$ws.Range("E38").Value2 = -0.009529860228716602
$ws.Range("D39").Value2 = 0.01134591630432031
$ws.Range("E39").Value2 = 0.005870342011230267
$ws.Range("D40").Value2 = 0.017133618377785
$ws.Range("E40").Value2 = -0.005201309959545375
$ws.Range("D41").Value2 = 0.0169168971587132
$ws.Range("E41").Value2 = 0.002797789746100676
$ws.Range("D42").Value2 = 0.03441943981861928
$ws.Range("E42").Value2 = 0.01126309504423828
$ws.Range("D43").Value2 = 0.01125742180653265
$ws.Range("E43").Value2 = 0.0009081353794406777
$ws.Range("D44").Value2 = 0.02220146971238447
$ws.Range("E44").Value2 = 0.001478260869565151
$ws.Range("D45").Value2 = 0.01306296489387153
$ws.Range("E45").Value2 = 0.01227125480016311
$ws.Range("D46").Value2 = 0.007783467851508971
$ws.Range("E46").Value2 = -0.007128969539857311
$ws.Range("D47").Value2 = 0.01289734130705501
$ws.Range("E47").Value2 = 0.001955586458650016
$ws.Range("D48").Value2 = 0.009550679930129825
$ws.Range("E48").Value2 = 0.01813380281690136
$ws.Range("D49").Value2 = 0.01519234428557056
$ws.Range("E49").Value2 = 0.006333239052104922
$ws.Range("D50").Value2 = 0.008773410523453049
$ws.Range("E50").Value2 = -0.02365150235308311
$ws.Range("D51").Value2 = 0.01162067895332743
$ws.Range("E51").Value2 = 0.02240097749720005
$ws.Range("D52").Value2 = 0.008763352915154171
$ws.Range("E52").Value2 = 0.006253664256400171
$ws.Range("D53").Value2 = 0.009778859972278989
$ws.Range("E53").Value2 = -0.005572396576319583
$ws.Range("D54").Value2 = 0.1352435066850867
$ws.Range("E54").Value2 = 0.0001970831690973895
$ws.Range("D55").Value2 = 0.04418466130038771
$ws.Range("E55").Value2 = 0.001677249311129891
$ws.Range("D56").Value2 = 0.9999999999999999
$ws.Range("E56").Value2 = 0.003994530778476557

$ws.Protect()
